# Commit: "replaced comma with semicolon in the
#   (<prog> -> program <identifier>; var <dec-list> begin <stat-list> end.)
#   rule in the grammar and parsing table"
#
# The document contains this rule twice: once in the grammar table and once
# in the parsing table. Each occurrence reads
#   "program <identifier>, var <dec-list> begin <stat-list> end."
# and the comma must become a semicolon. In the grammar table the result is
#   "program <identifier>; var <dec-list> begin <stat-list> end."
# (no extra space), while in the parsing table it is
#   "program <identifier> ; var <dec-list> begin <stat-list> end."
# (one extra space before the semicolon) -- matching how each table reads
# elsewhere in the document.
#
# Note: Find.Execute(..., Replace:=wdReplaceAll) in this host operates on
# the whole document regardless of which Range/Selection it was invoked on,
# so a scoped "replace all" would clobber unrelated "<identifier>, " text
# elsewhere (e.g. the "<identifier>, <dec>" rule). Instead we locate the
# exact character offsets of the target substring within the correct table
# cell and overwrite just that sub-range, leaving everything else intact.

$d = $word.ActiveDocument

$oldFragment = "<identifier>, "
$firstOccurrenceReplacement  = "<identifier>; "   # grammar table
$laterOccurrenceReplacement  = "<identifier> ; "  # parsing table(s)

$replacementsApplied = 0

for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $table = $d.Tables($t)
    $rowCount = $table.Rows.Count
    for ($row = 1; $row -le $rowCount; $row++) {
        $cell = $null
        try {
            $cell = $table.Cell($row, 2)
        } catch {
            $cell = $null
        }
        if ($cell -eq $null) {
            continue
        }

        $cellRange = $cell.Range
        $cellText = $cellRange.Text
        if ($cellText -eq $null) {
            continue
        }

        if ($cellText.StartsWith("program") -and $cellText.Contains($oldFragment)) {
            $idx = $cellText.IndexOf($oldFragment)
            $cellStart = $cellRange.Start
            $target = $d.Range($cellStart + $idx, $cellStart + $idx + $oldFragment.Length)

            if ($replacementsApplied -eq 0) {
                $target.Text = $firstOccurrenceReplacement
            } else {
                $target.Text = $laterOccurrenceReplacement
            }
            $replacementsApplied = $replacementsApplied + 1
        }
    }
}
